$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.411.56"
$ws.Range("E2").Value = "  +2.79%  "
$ws.Range("D3").Value = "3.390.83"
$ws.Range("E3").Value = "  +2.07%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.58"
$ws.Range("E5").Value = "  +2.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "180.75"
$ws.Range("E6").Value = "  +3.87%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("E8").Value = "  +1.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.196"
$ws.Range("E9").Value = "  +9.93%  "
$ws.Range("E10").Value = "  +1.95%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "48.98"
$ws.Range("E11").Value = "  +6.72%  "
$ws.Range("E12").Value = "  +5.15%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "686.97"
$ws.Range("E13").Value = "  -2.48%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.60"
$ws.Range("E14").Value = "  +2.22%  "
$ws.Range("D15").Value = "3.940.86"
$ws.Range("E15").Value = "  +2.01%  "
$ws.Range("D16").Value = "69.471.74"
$ws.Range("E16").Value = "  +2.87%  "
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.120"
$ws.Range("E17").Value = "  +1.80%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.380.96"
$ws.Range("E18").Value = "  +1.55%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.74"
$ws.Range("E19").Value = "  +2.26%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.41"
$ws.Range("E20").Value = "  +4.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.905"
$ws.Range("E21").Value = "  +1.65%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.39"
$ws.Range("E22").Value = "  +0.38%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.10"
$ws.Range("E23").Value = "  +1.29%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "104.46"
$ws.Range("E24").Value = "  +5.95%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.96"
$ws.Range("E25").Value = "  +2.87%  "
$ws.Range("E26").Value = "  +2.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.63"
$ws.Range("E27").Value = "  +3.29%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "34.72"
$ws.Range("E28").Value = "  +4.61%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.70"
$ws.Range("E29").Value = "  +2.45%  "
$ws.Range("E30").Value = "  -0.47%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.65"
$ws.Range("E32").Value = "  +9.78%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "555.84"
$ws.Range("E33").Value = "  -2.26%  "
$ws.Range("E34").Value = "  +1.76%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "58.43"
$ws.Range("E35").Value = "  +2.78%  "
$ws.Range("B36").Value = "Dai"
$ws.Range("C36").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.11%  "
$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D37").Value = "3.716.71"
$ws.Range("E37").Value = "  +0.50%  "
$ws.Range("E38").Value = "  +8.67%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "35.05"
$ws.Range("E39").Value = "  +2.61%  "
$ws.Range("D40").Value = "0.0₃0711"
$ws.Range("E40").Value = "  +6.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.26"
$ws.Range("E41").Value = "  +3.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.68"
$ws.Range("E42").Value = "  +2.26%  "
$ws.Range("E43").Value = "  +2.11%  "
$ws.Range("E44").Value = "  +3.65%  "
$ws.Range("E45").Value = "  -0.66%  "
$ws.Range("E46").Value = "  -1.26%  "
$ws.Range("E47").Value = "  +1.50%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.41"
$ws.Range("E48").Value = "  +7.30%  "
$ws.Range("E49").Value = "  +0.10%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "132.83"
$ws.Range("E50").Value = "  +2.92%  "
$ws.Range("E51").Value = "  -1.35%  "
